# Add week 33 ("semana 33 de 2025") column (AJ) to the weekly report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new week column, formatted like the other week headers
# (stored as text, same as the existing week-number headers, e.g. "32" in AI1).
# The leading apostrophe forces Excel to store the numeric-looking value as text.
$ws.Range("AJ1").Value = "'33"

# Weekly counts for week 33 per row (only rows that have a reported value).
$ws.Range("AJ2").Value = 56
$ws.Range("AJ5").Value = 3
$ws.Range("AJ6").Value = 71
$ws.Range("AJ7").Value = 20
$ws.Range("AJ8").Value = 12
$ws.Range("AJ9").Value = 1
$ws.Range("AJ10").Value = 3
$ws.Range("AJ11").Value = 1
$ws.Range("AJ13").Value = 1
$ws.Range("AJ14").Value = 3
$ws.Range("AJ15").Value = 1
$ws.Range("AJ23").Value = 3
$ws.Range("AJ24").Value = 1
$ws.Range("AJ25").Value = 58
$ws.Range("AJ26").Value = 4
$ws.Range("AJ29").Value = 0
$ws.Range("AJ31").Value = 6
$ws.Range("AJ34").Value = 0
$ws.Range("AJ35").Value = 25
$ws.Range("AJ36").Value = 5
$ws.Range("AJ37").Value = 6
$ws.Range("AJ38").Value = 88
$ws.Range("AJ40").Value = 0
$ws.Range("AJ41").Value = 4
$ws.Range("AJ42").Value = 7
$ws.Range("AJ43").Value = 24
$ws.Range("AJ45").Value = 60
$ws.Range("AJ46").Value = 122
$ws.Range("AJ47").Value = 0
$ws.Range("AJ48").Value = 68
$ws.Range("AJ49").Value = 4
$ws.Range("AJ50").Value = 0
$ws.Range("AJ51").Value = 13
$ws.Range("AJ53").Value = 5
$ws.Range("AJ54").Value = 0
$ws.Range("AJ55").Value = 0
$ws.Range("AJ56").Value = 8
$ws.Range("AJ58").Value = 30
